$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 94, pushing the existing rows 94-179 down to 95-180.
$ws.Rows.Item(94).Insert()

# Carry over the formatting/constant columns from the row above (same market /
# region / category / unit / origin / classification values shared by every
# row in this table) into the freshly inserted row.
$ws.Range("A93:R93").Copy()
$ws.Range("A94").PasteSpecial()

# Populate the new weekly record's own price data.
$ws.Range("D94").Value = 45159
$ws.Range("J94").Value = 250
$ws.Range("K94").Value = 11000
$ws.Range("L94").Value = 12000
$ws.Range("M94").Value = 11600
$ws.Range("P94").Value = 580
